# loginDetails.xlsx - "Added listener and created html report"
# A new login record (prakash@gmail.com / prakassh) is appended as row 5,
# mirroring the existing rows: column A is a mailto hyperlink styled like
# the other email cells, column B is the plain password text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "prakash@gmail.com"
$ws.Range("B5").Value = "prakassh"

# Link A5 to the new address and re-apply the same "Hyperlink" cell style
# that the existing linked cells (A2:A4) already use.
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:prakash@gmail.com")
$ws.Range("A5").Style = $ws.Range("A4").Style

# Move the selection the way it ended up after entering the new row.
$ws.Range("C5").Select() | Out-Null
